$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct / update individual plant.code values (case fixes & typo fixes) ---
$ws.Range("B2").Value   = "1u2"
$ws.Range("B3").Value   = "3e3"
$ws.Range("B4").Value   = "1u1"
$ws.Range("B25").Value  = "3e1"
$ws.Range("B81").Value  = "10e9"
$ws.Range("B86").Value  = "9u1"
$ws.Range("B99").Value  = "10e6"
$ws.Range("B106").Value = "10u8"
$ws.Range("B107").Value = "9u9"
$ws.Range("B132").Value = "10e1"

# --- Row 136 no longer flagged as "plant.code.unsure" ---
$ws.Range("D136").ClearContents()

# --- Turn on AutoFilter for the data range ---
$ws.Range("A1:E164").AutoFilter() | Out-Null

# --- Register the (hidden) filter-database defined name tied to the sheet ---
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$164")
$fdb.Visible = $false

# --- Update the view: selection + zoom ---
$ws.Range("E5").Select()
$excel.ActiveWindow.Zoom = 150
